# edit.ps1 - PowerPoint COM-interop script
#
# Reproduces the authoring session captured by the diff:
#   - The footer "Date" placeholder field (datetimeFigureOut) is refreshed
#     from 7/23/2020 -> 8/20/2020 on the slide master and on every custom
#     (slide) layout that carries a date placeholder.
#   - On slide 1 ("Tableau 4" table), the row-5 / column-1 cell text is
#     edited from "Input Metadata (Optional)" to "Input Metadata (Suggested)".
#
# ppPlaceholderDate = 16 ; msoPlaceholder = 14
$ppPlaceholderDate = 16
$msoPlaceholder = 14
$NEW_DATE = "8/20/2020"

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Type -eq $msoPlaceholder) {
            $pf = $shp.PlaceholderFormat
            if ($pf.Type -eq $ppPlaceholderDate) {
                $shp.TextFrame.TextRange.Text = $NEW_DATE
            }
        }
    }
}

$p = $ppt.ActivePresentation

# --- Refresh the cached "datetimeFigureOut" footer field -------------------
# Slide master
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every custom (slide) layout off the master
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i).Shapes
}

# --- Slide 1: update the table cell text ------------------------------------
$s = $p.Slides.Item(1)
$shp = $s.Shapes.Item(1)
if ($shp.HasTable) {
    $tbl = $shp.Table
    $cell = $tbl.Cell(5, 1)
    $cell.Shape.TextFrame.TextRange.Text = "Input Metadata (Suggested)"
}
